# Generate Report for Handoff
# Updates the status/dates/error-detail for the f3e2fe82-... file row
# across the Overview, zh-cn and de-de worksheets, reflecting that a new
# handoff round has started (status back to "Ready for handoff") and that
# the previous handback is stale (Error Detail populated).

$wb = $excel.ActiveWorkbook

$newHandoffDateTime = "2016-09-02 18:54:22"
$zhHandoffDateTime   = "2016-09-02 18:54:15"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bc3d193b46e5305afee82e013ed50d654ddbf568/e2e/f3e2fe82-073d-4902-8608-044de7c6793f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c15bfcb22ea774523fe631df3bc20ccf64c97f30/e2e/f3e2fe82-073d-4902-8608-044de7c6793f.md."

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = $newHandoffDateTime

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("H3").Value = $zhHandoffDateTime
$wsZh.Range("P3").Value = $errorDetail
# 39.1 (Excel "characters" units) is the value that round-trips to an
# on-disk column width of exactly 40 once Excel converts characters -> pixels
# -> characters internally (the naive ColumnWidth = 40 serializes as ~40.83).
$wsZh.Columns.Item(16).ColumnWidth = 39.1

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("H3").Value = $newHandoffDateTime
$wsDe.Range("P3").Value = $errorDetail
$wsDe.Columns.Item(16).ColumnWidth = 39.1
